# Actualización automática 2025-08-04 17:26:10
# Add a new sheet "CUMPLIMIENTO MENSUAL" after "VENTA MENSUAL" with the
# monthly compliance summary for VACA PANCHI CAROLINA.

$wb = $excel.ActiveWorkbook

$ventaMensual = $wb.Worksheets.Item("VENTA MENSUAL")
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ventaMensual)
$ws.Name = "CUMPLIMIENTO MENSUAL"

# --- column widths (match VENTAS POR GRUPO / VENTA MENSUAL style sheets) ---
$offset = 5.0 / 6.0
$ws.Columns.Item(1).ColumnWidth = 22 - $offset
$ws.Columns.Item(2).ColumnWidth = 13 - $offset
$ws.Columns.Item(3).ColumnWidth = 17 - $offset
$ws.Columns.Item(4).ColumnWidth = 11 - $offset
$ws.Columns.Item(5).ColumnWidth = 17 - $offset
$ws.Columns.Item(6).ColumnWidth = 18 - $offset

# --- header row (reuse the existing bold/border/center-top header style) ---
$headerSrc = $ventaMensual.Range("A1")
$headerSrc.Copy()
$ws.Range("A1:F1").PasteSpecial(-4122)

$ws.Range("A1").Value = "ASESOR"
$ws.Range("B1").Value = "GRUPO"
$ws.Range("C1").Value = "PRESUPUESTO"
$ws.Range("D1").Value = "VENTA"
$ws.Range("E1").Value = "POR CUMPLIR"
$ws.Range("F1").Value = "CUMPLIMIENTO"

# --- data rows ---
$ws.Range("A2").Value = "VACA PANCHI CAROLINA"
$ws.Range("B2").Value = "OTROS"
$ws.Range("C2").Value = 0
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = 0
$ws.Range("F2").Value = 0

$ws.Range("A3").Value = "VACA PANCHI CAROLINA"
$ws.Range("B3").Value = "PORCELANATO"
$ws.Range("C3").Value = 10000
$ws.Range("D3").Value = 0
$ws.Range("E3").Value = 10000
$ws.Range("F3").Value = 0

$ws.Range("B4").Value = "TOTAL"
$ws.Range("C4").Value = 10000
$ws.Range("D4").Value = 0
$ws.Range("E4").Value = 10000
$ws.Range("F4").Value = 0

# --- number formats for PRESUPUESTO / VENTA / POR CUMPLIR (left aligned $) ---
$ws.Range("C2:E4").NumberFormat = """$""#,##0.00"

# --- number format + right alignment for CUMPLIMIENTO column ---
$ws.Range("F2:F4").NumberFormat = """$""#,##0.00"
$ws.Range("F2:F4").HorizontalAlignment = -4152

# --- TOTAL label right aligned (new style: general format + right align) ---
$ws.Range("B4").HorizontalAlignment = -4152

# Keep the first sheet active, like in the original workbook.
$wb.Worksheets.Item(1).Activate()
